# The workbook is a set of per-teacher weekly schedule sheets. Cells in the
# "grupa" (group) columns that have no real group assigned use the literal
# placeholder text "X". This edit turns every such placeholder into "-"
# across every worksheet in the workbook.

$wb = $excel.ActiveWorkbook

$replaced = 0

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            if ($cell.Value() -eq "X") {
                $cell.Value = "-"
                $replaced = $replaced + 1
            }
        }
    }
}

Write-Host "Replaced $replaced cell(s) containing 'X' with '-'"
